# Update the test email addresses in the TestData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")

$ws.Range("B2").Value = "test26@mail.com"
$ws.Range("B3").Value = "test27@mail.com"
